$wb = $excel.ActiveWorkbook

# Both "展览" (sheet 1) and "全部类型" (sheet 4) carry identical data rows
# and both receive the same F (想去人数) / G (最低票价) value updates.
$targetSheets = @("展览", "全部类型")

# Map of cell reference -> new value, taken from the authoritative diff.
$updates = @{
    "F2" = 1187
    "G2" = 65
    "F3" = 965
    "F4" = 293
    "G4" = 50
    "F5" = 61
    "G5" = 40
    "F6" = 1124
    "G6" = 65
    "G7" = 158
    "F8" = 2424
    "F9" = 7911
    "G9" = 75
    "F10" = 941
    "G10" = 60
    "F11" = 473
    "G11" = 68
    "F12" = 416
    "G12" = 65
    "F13" = 179
    "G13" = 25
    "F14" = 445
    "F15" = 9
    "F16" = 170
    "F17" = 8164
    "F19" = 1411
    "F20" = 165
    "F23" = 191
    "F24" = 344
    "F25" = 191
    "F29" = 36
    "F30" = 433
    "F31" = 1169
    "F32" = 27
    "F34" = 104
    "F36" = 89
}

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($ref in $updates.Keys) {
        $ws.Range($ref).Value = $updates[$ref]
    }
}

Write-Output "Updated $($updates.Count) cells on each of: $($targetSheets -join ', ')"
